$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.592.05'
$ws.Range("E2").Value = '  -2.21%  '
$ws.Range("D3").Value = '2.888.23'
$ws.Range("E3").Value = '  -2.69%  '
$ws.Range("D4").Value = '0.994'
$ws.Range("E4").Value = '  -0.54%  '
$ws.Range("D5").Value = '572.32'
$ws.Range("E5").Value = '  -3.95%  '
$ws.Range("D6").Value = '143.55'
$ws.Range("E6").Value = '  -2.40%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '0.501'
$ws.Range("E8").Value = '  -1.30%  '
$ws.Range("D9").Value = '2.876.93'
$ws.Range("E9").Value = '  -3.01%  '
$ws.Range("D10").Value = '6.71'
$ws.Range("E10").Value = '  -7.38%  '
$ws.Range("D11").Value = '0.148'
$ws.Range("E11").Value = '  -2.99%  '
$ws.Range("D12").Value = '0.431'
$ws.Range("E12").Value = '  -3.41%  '
$ws.Range("D13").Value = '0.0000233'
$ws.Range("E13").Value = '  -3.38%  '
$ws.Range("D14").Value = '32.08'
$ws.Range("E14").Value = '  -4.00%  '
$ws.Range("D15").Value = '0.125'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").Value = '3.354.31'
$ws.Range("E16").Value = '  -3.15%  '
$ws.Range("D17").Value = '61.332.32'
$ws.Range("E17").Value = '  -2.45%  '
$ws.Range("D18").Value = '6.59'
$ws.Range("E18").Value = '  -2.28%  '
$ws.Range("D19").Value = '2.890.11'
$ws.Range("E19").Value = '  -2.79%  '
$ws.Range("D20").Value = '434.14'
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("D21").Value = '13.20'
$ws.Range("E21").Value = '  -2.66%  '
$ws.Range("D22").Value = '0.654'
$ws.Range("E22").Value = '  -2.95%  '
$ws.Range("D23").Value = '6.89'
$ws.Range("E23").Value = '  -3.04%  '
$ws.Range("D24").Value = '79.20'
$ws.Range("E24").Value = '  -3.04%  '
$ws.Range("D25").Value = '11.85'
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D26").Value = '10.08'
$ws.Range("E26").Value = '  -10.59%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").Value = '2.03'
$ws.Range("E28").Value = '  -5.25%  '
$ws.Range("D29").Value = '0.0000108'
$ws.Range("E29").Value = '  +11.11%  '
$ws.Range("D30").Value = '7.00'
$ws.Range("E30").Value = '  -3.62%  '
$ws.Range("D31").Value = '2.51'
$ws.Range("E31").Value = '  -4.46%  '
$ws.Range("D32").Value = '2.07'
$ws.Range("E32").Value = '  -4.63%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("D34").Value = '0.106'
$ws.Range("E34").Value = '  -3.86%  '
$ws.Range("D35").Value = '25.59'
$ws.Range("E35").Value = '  -3.91%  '
$ws.Range("D36").Value = '0.955'
$ws.Range("E36").Value = '  -4.10%  '
$ws.Range("D37").Value = '5.43'
$ws.Range("E37").Value = '  -4.32%  '
$ws.Range("D38").Value = '2.97'
$ws.Range("E38").Value = '  -4.45%  '
$ws.Range("D39").Value = '49.02'
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("D40").Value = '1.95'
$ws.Range("E40").Value = '  -4.74%  '
$ws.Range("D41").Value = '0.116'
$ws.Range("E41").Value = '  -2.49%  '
$ws.Range("D42").Value = '8.25'
$ws.Range("E42").Value = '  -3.51%  '
$ws.Range("D43").Value = '0.268'
$ws.Range("E43").Value = '  -5.13%  '
$ws.Range("D44").Value = '38.23'
$ws.Range("E44").Value = '  -5.93%  '
$ws.Range("D45").Value = '2.679.06'
$ws.Range("E45").Value = '  -2.42%  '
$ws.Range("D46").Value = '133.01'
$ws.Range("E46").Value = '  -1.34%  '
$ws.Range("D47").Value = '0.0331'
$ws.Range("E47").Value = '  -2.88%  '
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").Value = '336.46'
$ws.Range("E49").Value = '  -8.22%  '
$ws.Range("D50").Value = '0.102'
$ws.Range("E50").Value = '  -2.67%  '
$ws.Range("D51").Value = '21.65'
$ws.Range("E51").Value = '  -6.24%  '
